$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.828516
$ws.Range("H2").Value = 2.485548
$ws.Range("I2").Value = 0.4625620436231038
$ws.Range("J2").Value = 0.4821955800271095
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.472365
$ws.Range("N2").Value = 13.417095
$ws.Range("O2").Value = 0.1840876942178652
$ws.Range("P2").Value = 0.2075728609309428
$ws.Range("Q2").Value = 3.70542596034
$ws.Range("R2").Value = 33.34883364306
$ws.Range("S2").Value = 0.08515198004328077
$ws.Range("T2").Value = 0.1000907160744825
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.828516
$ws.Range("H3").Value = 2.485548
$ws.Range("I3").Value = 0.4625620436231038
$ws.Range("J3").Value = 0.4821955800271095
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.966196333333333
$ws.Range("N3").Value = 11.898589
$ws.Range("O3").Value = 0.163253208943967
$ws.Range("P3").Value = 0.1840803959256042
$ws.Range("Q3").Value = 3.286057121308
$ws.Range("R3").Value = 29.574514091772
$ws.Range("S3").Value = 0.07551473795715094
$ws.Range("T3").Value = 0.08876275328496672
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.828516
$ws.Range("H4").Value = 2.485548
$ws.Range("I4").Value = 0.4625620436231038
$ws.Range("J4").Value = 0.4821955800271095
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.115150666666667
$ws.Range("N4").Value = 12.345452
$ws.Range("O4").Value = 0.1693843408545093
$ws.Range("P4").Value = 0.1909937129554221
$ws.Range("Q4").Value = 3.409468169744
$ws.Range("R4").Value = 30.685213527696
$ws.Range("S4").Value = 0.07835076686341423
$ws.Range("T4").Value = 0.09209632420007104
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.828516
$ws.Range("H5").Value = 2.485548
$ws.Range("I5").Value = 0.4625620436231038
$ws.Range("J5").Value = 0.4821955800271095
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.494784666666666
$ws.Range("N5").Value = 10.484354
$ws.Range("O5").Value = 0.1438493618196675
$ws.Range("P5").Value = 0.1622010841238564
$ws.Range("Q5").Value = 2.895485012888
$ws.Range("R5").Value = 26.059365115992
$ws.Range("S5").Value = 0.0665392547771847
$ws.Range("T5").Value = 0.07821264584012894
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.828516
$ws.Range("H6").Value = 2.485548
$ws.Range("I6").Value = 0.4625620436231038
$ws.Range("J6").Value = 0.4821955800271095
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 8.246256000000001
$ws.Range("N6").Value = 16.492512
$ws.Range("O6").Value = 0.3394253941639908
$ws.Range("P6").Value = 0.2551519460641745
$ws.Range("Q6").Value = 6.832155036096001
$ws.Range("R6").Value = 40.99293021657601
$ws.Range("S6").Value = 0.1570053039820731
$ws.Range("T6").Value = 0.1230331406274604
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7438396666666667
$ws.Range("H7").Value = 2.231519
$ws.Range("I7").Value = 0.4152870872032183
$ws.Range("J7").Value = 0.4329140288365043
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.472365
$ws.Range("N7").Value = 13.417095
$ws.Range("O7").Value = 0.1840876942178652
$ws.Range("P7").Value = 0.2075728609309428
$ws.Range("Q7").Value = 3.326722490811667
$ws.Range("R7").Value = 29.940502417305
$ws.Range("S7").Value = 0.07644924232169399
$ws.Range("T7").Value = 0.08986120350273384
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7438396666666667
$ws.Range("H8").Value = 2.231519
$ws.Range("I8").Value = 0.4152870872032183
$ws.Range("J8").Value = 0.4329140288365043
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.966196333333333
$ws.Range("N8").Value = 11.898589
$ws.Range("O8").Value = 0.163253208943967
$ws.Range("P8").Value = 0.1840803959256042
$ws.Range("Q8").Value = 2.950214158521222
$ws.Range("R8").Value = 26.551927426691
$ws.Range("S8").Value = 0.06779694961891844
$ws.Range("T8").Value = 0.07969098582997215
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7438396666666667
$ws.Range("H9").Value = 2.231519
$ws.Range("I9").Value = 0.4152870872032183
$ws.Range("J9").Value = 0.4329140288365043
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.115150666666667
$ws.Range("N9").Value = 12.345452
$ws.Range("O9").Value = 0.1693843408545093
$ws.Range("P9").Value = 0.1909937129554221
$ws.Range("Q9").Value = 3.061012300176444
$ws.Range("R9").Value = 27.549110701588
$ws.Range("S9").Value = 0.07034312953130628
$ws.Range("T9").Value = 0.08268385775797463
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7438396666666667
$ws.Range("H10").Value = 2.231519
$ws.Range("I10").Value = 0.4152870872032183
$ws.Range("J10").Value = 0.4329140288365043
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.494784666666666
$ws.Range("N10").Value = 10.484354
$ws.Range("O10").Value = 0.1438493618196675
$ws.Range("P10").Value = 0.1622010841238564
$ws.Range("Q10").Value = 2.599559461525111
$ws.Range("R10").Value = 23.396035153726
$ws.Range("S10").Value = 0.05973878246613157
$ws.Range("T10").Value = 0.07021912480970743
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.7438396666666667
$ws.Range("H11").Value = 2.231519
$ws.Range("I11").Value = 0.4152870872032183
$ws.Range("J11").Value = 0.4329140288365043
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 8.246256000000001
$ws.Range("N11").Value = 16.492512
$ws.Range("O11").Value = 0.3394253941639908
$ws.Range("P11").Value = 0.2551519460641745
$ws.Range("Q11").Value = 6.133892314288
$ws.Range("R11").Value = 36.803353885728
$ws.Range("S11").Value = 0.140958983265168
$ws.Range("T11").Value = 0.1104588569361162
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.21879
$ws.Range("H12").Value = 0.43758
$ws.Range("I12").Value = 0.1221508691736778
$ws.Range("J12").Value = 0.08489039113638626
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.472365
$ws.Range("N12").Value = 13.417095
$ws.Range("O12").Value = 0.1840876942178652
$ws.Range("P12").Value = 0.2075728609309428
$ws.Range("Q12").Value = 0.97850873835
$ws.Range("R12").Value = 5.8710524301
$ws.Range("S12").Value = 0.02248647185289047
$ws.Range("T12").Value = 0.01762094135372644
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.21879
$ws.Range("H13").Value = 0.43758
$ws.Range("I13").Value = 0.1221508691736778
$ws.Range("J13").Value = 0.08489039113638626
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.966196333333333
$ws.Range("N13").Value = 11.898589
$ws.Range("O13").Value = 0.163253208943967
$ws.Range("P13").Value = 0.1840803959256042
$ws.Range("Q13").Value = 0.8677640957699999
$ws.Range("R13").Value = 5.206584574619999
$ws.Range("S13").Value = 0.01994152136789761
$ws.Range("T13").Value = 0.01562665681066539
# Row 14
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.21879
$ws.Range("H14").Value = 0.43758
$ws.Range("I14").Value = 0.1221508691736778
$ws.Range("J14").Value = 0.08489039113638626
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.115150666666667
$ws.Range("N14").Value = 12.345452
$ws.Range("O14").Value = 0.1693843408545093
$ws.Range("P14").Value = 0.1909937129554221
$ws.Range("Q14").Value = 0.9003538143599999
$ws.Range("R14").Value = 5.40212288616
$ws.Range("S14").Value = 0.02069044445978883
$ws.Range("T14").Value = 0.01621353099737647
# Row 15
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.21879
$ws.Range("H15").Value = 0.43758
$ws.Range("I15").Value = 0.1221508691736778
$ws.Range("J15").Value = 0.08489039113638626
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.494784666666666
$ws.Range("N15").Value = 10.484354
$ws.Range("O15").Value = 0.1438493618196675
$ws.Range("P15").Value = 0.1622010841238564
$ws.Range("Q15").Value = 0.7646239372199999
$ws.Range("R15").Value = 4.58774362332
$ws.Range("S15").Value = 0.01757132457635126
$ws.Range("T15").Value = 0.01376931347402006
# Row 16
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.21879
$ws.Range("H16").Value = 0.43758
$ws.Range("I16").Value = 0.1221508691736778
$ws.Range("J16").Value = 0.08489039113638626
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 8.246256000000001
$ws.Range("N16").Value = 16.492512
$ws.Range("O16").Value = 0.3394253941639908
$ws.Range("P16").Value = 0.2551519460641745
$ws.Range("Q16").Value = 1.80419835024
$ws.Range("R16").Value = 7.21679340096
$ws.Range("S16").Value = 0.04146110691674968
$ws.Range("T16").Value = 0.0216599485005979
